$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44694
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 16000
$ws.Range("L3").Value = 17000
$ws.Range("M3").Value = 16500
$ws.Range("P3").Value = 1650

$ws.Range("D4").Value = 44679
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 19000
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = 19500
$ws.Range("P4").Value = 1950

$ws.Range("D5").Value = 44406
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14500
$ws.Range("P5").Value = 1450

$ws.Range("D6").Value = 44727
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 18000
$ws.Range("L6").Value = 19000
$ws.Range("M6").Value = 18500
$ws.Range("P6").Value = 1850

$ws.Range("D7").Value = 44160
$ws.Range("J7").Value = 360
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 11000
$ws.Range("M7").Value = 10500
$ws.Range("P7").Value = 1050

$ws.Range("D8").Value = 44714
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 19000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 19500
$ws.Range("P8").Value = 1950

$ws.Range("D9").Value = 44580
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 19000
$ws.Range("P9").Value = 1900

$ws.Range("D10").Value = 44524
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 20000
$ws.Range("L10").Value = 21000
$ws.Range("M10").Value = 20500
$ws.Range("P10").Value = 2050

$ws.Range("D11").Value = 44644
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 20000
$ws.Range("L11").Value = 21000
$ws.Range("M11").Value = 20500
$ws.Range("P11").Value = 2050

$ws.Range("D12").Value = 44291
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 13000
$ws.Range("L12").Value = 14000
$ws.Range("M12").Value = 13500
$ws.Range("P12").Value = 1350

$ws.Range("D13").Value = 44265
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 15500
$ws.Range("P13").Value = 1550

$ws.Range("D14").Value = 44428
$ws.Range("J14").Value = 300
$ws.Range("K14").Value = 15000
$ws.Range("L14").Value = 16000
$ws.Range("M14").Value = 15500
$ws.Range("P14").Value = 1550

$ws.Range("D15").Value = 44263
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 15000
$ws.Range("L15").Value = 16000
$ws.Range("M15").Value = 15500
$ws.Range("P15").Value = 1550

$ws.Range("D16").Value = 44218
$ws.Range("J16").Value = 320
$ws.Range("K16").Value = 10000
$ws.Range("L16").Value = 11000
$ws.Range("M16").Value = 10500
$ws.Range("P16").Value = 1050

$ws.Range("D17").Value = 44460
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 16000
$ws.Range("M17").Value = 15500
$ws.Range("P17").Value = 1550

$ws.Range("D18").Value = 44441
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 15000
$ws.Range("L18").Value = 16000
$ws.Range("M18").Value = 15500
$ws.Range("P18").Value = 1550

$ws.Range("D19").Value = 44204
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 10000
$ws.Range("L19").Value = 11000
$ws.Range("M19").Value = 10500
$ws.Range("P19").Value = 1050

$ws.Range("D20").Value = 44377
$ws.Range("J20").Value = 650
$ws.Range("K20").Value = 14000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 14538
$ws.Range("P20").Value = 1454

$ws.Range("D21").Value = 44547
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 19000
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = 19500
$ws.Range("P21").Value = 1950

$ws.Range("D22").Value = 44358
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 14000
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = 14500
$ws.Range("P22").Value = 1450
